$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.914.39"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "3.803.23"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "699.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("D7").Value = "3.801.22"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("E13").Value = "  +9.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").Value = "4.448.92"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "3.805.33"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "70.932.16"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "3.954.97"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +15.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.95%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "3.754.81"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.43%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.972"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +23.79%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "162.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.300"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("E51").Value = "  -1.23%  "
